# category-fiction.xlsx edit: add STT (row number), TITTLE (book title),
# ISBN13 and Cost columns to the Sheet1 book table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Make room: insert two columns before the current column A
#    (Author's Name, Publisher, Language, Country, Publication Year,
#    Description -> shift from A:F to C:H), then insert two more columns
#    before the (now shifted) Country column F, so ISBN13/Cost land
#    between Language and Country.
# ---------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert()
$ws.Range("F1:G1").EntireColumn.Insert()

# Final layout:
# A STT | B TITTLE | C Author's Name | D Publisher | E Language |
# F ISBN13 | G Cost | H Country | I Publication Year | J Description

# ---------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "STT"
$ws.Range("B1").Value = "TITTLE"
$ws.Range("F1").Value = "ISBN13"
$ws.Range("G1").Value = "Cost"

$ws.Range("A1:J1").Font.Bold = $true
$ws.Range("F1").NumberFormat = "0"

# ---------------------------------------------------------------------
# 3. STT column (A2:A11) - simple row numbers
# ---------------------------------------------------------------------
$sttValues = @(1,2,3,4,5,6,7,8,9,10)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $sttValues[$i]
}

# ---------------------------------------------------------------------
# 4. TITTLE column (B2:B11)
# ---------------------------------------------------------------------
$titles = @(
    "The Girl on the Train",
    "The Sellout",
    "The Tales of Beedle the Bard",
    "DanTDM: Trayaurus and the Enchanted Crystal`n",
    "Night School",
    "Harry Potter and the Philosopher's Stone",
    "The Agent",
    "The Last Paradise",
    "A Criminal Defense",
    "Executive Order"
)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $titles[$i]
}

# ---------------------------------------------------------------------
# 5. ISBN13 column (F2:F11) and Cost column (G2:G11)
# ---------------------------------------------------------------------
$isbns = @(
    9780552779777,
    9781786070159,
    9780747599876,
    9781409168393,
    9780593073919,
    9781408845646,
    9781477818022,
    9781503941885,
    9781503943421,
    9781477819432
)
$costs = @(
    9.5399999999999991,
    12.66,
    8.5500000000000007,
    17.13,
    17.600000000000001,
    41.2,
    13.82,
    14.34,
    13.46,
    12.79
)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $isbns[$i]
    $ws.Cells.Item($i + 2, 7).Value = $costs[$i]
}

$ws.Range("F2:F11").NumberFormat = "0"
$ws.Range("F2:F11").Font.Color = 0x333333
$ws.Range("G2:G11").NumberFormat = """$""#,##0.00_);[Red](""$""#,##0.00)"
$ws.Range("G2:G11").Font.Color = 0x7200FF

# ---------------------------------------------------------------------
# 6. TITTLE column formatting
#    - B2: #333333 colour, wrap text, vertically centred (not bold)
#    - B3:B11 (except B5): bold, #333333 colour, wrap text, vertically
#      centred
#    - B5: wrap text only, default colour/weight, taller row
# ---------------------------------------------------------------------
$ws.Range("B2").Font.Color = 0x333333
$ws.Range("B2").WrapText = $true
$ws.Range("B2").VerticalAlignment = -4108

$ws.Range("B3:B4").Font.Bold = $true
$ws.Range("B3:B4").Font.Color = 0x333333
$ws.Range("B3:B4").WrapText = $true
$ws.Range("B3:B4").VerticalAlignment = -4108

$ws.Range("B6:B11").Font.Bold = $true
$ws.Range("B6:B11").Font.Color = 0x333333
$ws.Range("B6:B11").WrapText = $true
$ws.Range("B6:B11").VerticalAlignment = -4108

$ws.Range("B5").WrapText = $true

# ---------------------------------------------------------------------
# 7. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(5).RowHeight = 60

# ---------------------------------------------------------------------
# 8. Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 38
$ws.Columns.Item(3).ColumnWidth = 28.5703125
$ws.Columns.Item(4).ColumnWidth = 30.85546875
$ws.Columns.Item(5).ColumnWidth = 10
$ws.Columns.Item(6).ColumnWidth = 14.140625
$ws.Columns.Item(7).ColumnWidth = 15.5703125
$ws.Columns.Item(8).ColumnWidth = 20.28515625
$ws.Columns.Item(9).ColumnWidth = 15.7109375
$ws.Columns.Item(10).ColumnWidth = 255.7109375

# ---------------------------------------------------------------------
# 9. View / selection tidy-up
# ---------------------------------------------------------------------
$ws.Range("F1:F1048576").Select()
